{"js": "// Update the date paragraph (first paragraph in the body, outside the table)\nconst dateParagraph = context.document.body.paragraphs.getFirst();\ndateParagraph.insertText(\"2025-04-06 Sunday\", Word.InsertLocation.replace);\n\n// Update the table of arithmetic answers.\n// table.values preserves each cell's existing paragraph/run formatting\n// while replacing the text content.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = [\n  [\n    \"29+66=95\",\n    \"68+24=92\",\n    \"8+57=65\",\n    \"17+15=32\",\n    \"5+39=44\"\n  ],\n  [\n    \"42-24=18\",\n    \"59+19=78\",\n    \"6+28=34\",\n    \"20-19=1\",\n    \"83-38=45\"\n  ],\n  [\n    \"70-64=6\",\n    \"9+17=26\",\n    \"45+17=62\",\n    \"71-2=69\",\n    \"51-34=17\"\n  ],\n  [\n    \"34-29=5\",\n    \"53-14=39\",\n    \"43-15=28\",\n    \"83-6=77\",\n    \"56+25=81\"\n  ],\n  [\n    \"62-36=26\",\n    \"81-44=37\",\n    \"44-6=38\",\n    \"4+57=61\",\n    \"9+56=65\"\n  ],\n  [\n    \"97-89=8\",\n    \"6+58=64\",\n    \"17+57=74\",\n    \"47+16=63\",\n    \"58+26=84\"\n  ],\n  [\n    \"17-9=8\",\n    \"5+48=53\",\n    \"66+8=74\",\n    \"46-8=38\",\n    \"70-56=14\"\n  ],\n  [\n    \"69+3=72\",\n    \"39+23=62\",\n    \"81-56=25\",\n    \"69+19=88\",\n    \"61-52=9\"\n  ],\n  [\n    \"73-15=58\",\n    \"57+37=94\",\n    \"46-18=28\",\n    \"17+57=74\",\n    \"53-14=39\"\n  ],\n  [\n    \"17+78=95\",\n    \"47+25=72\",\n    \"29+43=72\",\n    \"52-25=27\",\n    \"92-4=88\"\n  ],\n  [\n    \"72-17=55\",\n    \"61-14=47\",\n    \"61-13=48\",\n    \"47+15=62\",\n    \"69+7=76\"\n  ],\n  [\n    \"81-77=4\",\n    \"38+46=84\",\n    \"64-25=39\",\n    \"76-47=29\",\n    \"96-9=87\"\n  ],\n  [\n    \"63-6=57\",\n    \"90-71=19\",\n    \"17+27=44\",\n    \"46+48=94\",\n    \"75+6=81\"\n  ],\n  [\n    \"63-49=14\",\n    \"17+69=86\",\n    \"73-25=48\",\n    \"41-33=8\",\n    \"13+48=61\"\n  ],\n  [\n    \"28+47=75\",\n    \"39+44=83\",\n    \"33+48=81\",\n    \"18+55=73\",\n    \"72-8=64\"\n  ],\n  [\n    \"18+45=63\",\n    \"6+88=94\",\n    \"36-9=27\",\n    \"77+14=91\",\n    \"59+27=86\"\n  ],\n  [\n    \"95-56=39\",\n    \"30-12=18\",\n    \"26+9=35\",\n    \"53-44=9\",\n    \"19+64=83\"\n  ],\n  [\n    \"43-29=14\",\n    \"37-18=19\",\n    \"8+75=83\",\n    \"70-22=48\",\n    \"91-72=19\"\n  ],\n  [\n    \"38+6=44\",\n    \"40-4=36\",\n    \"7+67=74\",\n    \"17+5=22\",\n    \"19+19=38\"\n  ],\n  [\n    \"35+8=43\",\n    \"93-26=67\",\n    \"69+19=88\",\n    \"83-45=38\",\n    \"9+53=62\"\n  ]\n];\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date paragraph (first paragraph in the body, outside the table)\n$d.Paragraphs(1).Range.Text = \"2025-04-06 Sunday\"\n\n# New values for the 100 table cells, in row-major (reading) order:\n# row 1 col 1, row 1 col 2, ..., row 1 col 5, row 2 col 1, ...\n$newValues = @(\n    \"29+66=95\",\n    \"68+24=92\",\n    \"8+57=65\",\n    \"17+15=32\",\n    \"5+39=44\",\n    \"42-24=18\",\n    \"59+19=78\",\n    \"6+28=34\",\n    \"20-19=1\",\n    \"83-38=45\",\n    \"70-64=6\",\n    \"9+17=26\",\n    \"45+17=62\",\n    \"71-2=69\",\n    \"51-34=17\",\n    \"34-29=5\",\n    \"53-14=39\",\n    \"43-15=28\",\n    \"83-6=77\",\n    \"56+25=81\",\n    \"62-36=26\",\n    \"81-44=37\",\n    \"44-6=38\",\n    \"4+57=61\",\n    \"9+56=65\",\n    \"97-89=8\",\n    \"6+58=64\",\n    \"17+57=74\",\n    \"47+16=63\",\n    \"58+26=84\",\n    \"17-9=8\",\n    \"5+48=53\",\n    \"66+8=74\",\n    \"46-8=38\",\n    \"70-56=14\",\n    \"69+3=72\",\n    \"39+23=62\",\n    \"81-56=25\",\n    \"69+19=88\",\n    \"61-52=9\",\n    \"73-15=58\",\n    \"57+37=94\",\n    \"46-18=28\",\n    \"17+57=74\",\n    \"53-14=39\",\n    \"17+78=95\",\n    \"47+25=72\",\n    \"29+43=72\",\n    \"52-25=27\",\n    \"92-4=88\",\n    \"72-17=55\",\n    \"61-14=47\",\n    \"61-13=48\",\n    \"47+15=62\",\n    \"69+7=76\",\n    \"81-77=4\",\n    \"38+46=84\",\n    \"64-25=39\",\n    \"76-47=29\",\n    \"96-9=87\",\n    \"63-6=57\",\n    \"90-71=19\",\n    \"17+27=44\",\n    \"46+48=94\",\n    \"75+6=81\",\n    \"63-49=14\",\n    \"17+69=86\",\n    \"73-25=48\",\n    \"41-33=8\",\n    \"13+48=61\",\n    \"28+47=75\",\n    \"39+44=83\",\n    \"33+48=81\",\n    \"18+55=73\",\n    \"72-8=64\",\n    \"18+45=63\",\n    \"6+88=94\",\n    \"36-9=27\",\n    \"77+14=91\",\n    \"59+27=86\",\n    \"95-56=39\",\n    \"30-12=18\",\n    \"26+9=35\",\n    \"53-44=9\",\n    \"19+64=83\",\n    \"43-29=14\",\n    \"37-18=19\",\n    \"8+75=83\",\n    \"70-22=48\",\n    \"91-72=19\",\n    \"38+6=44\",\n    \"40-4=36\",\n    \"7+67=74\",\n    \"17+5=22\",\n    \"19+19=38\",\n    \"35+8=43\",\n    \"93-26=67\",\n    \"69+19=88\",\n    \"83-45=38\",\n    \"9+53=62\"\n)\n\n$t = $d.Tables(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$i]\n        $i++\n    }\n}\n"}
